$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-08-17 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-18 Sunday", 2)

$d.Content.Find.Execute("13×45=", $true, $false, $false, $false, $false, $true, 1, $false, "14×25=", 2)
$d.Content.Find.Execute("41×66=", $true, $false, $false, $false, $false, $true, 1, $false, "93×31=", 2)
$d.Content.Find.Execute("35×63=", $true, $false, $false, $false, $false, $true, 1, $false, "88×76=", 2)
$d.Content.Find.Execute("82×45=", $true, $false, $false, $false, $false, $true, 1, $false, "70×33=", 2)
$d.Content.Find.Execute("39×81=", $true, $false, $false, $false, $false, $true, 1, $false, "69×54=", 2)
$d.Content.Find.Execute("38×11=", $true, $false, $false, $false, $false, $true, 1, $false, "36×14=", 2)
$d.Content.Find.Execute("53×65=", $true, $false, $false, $false, $false, $true, 1, $false, "13×79=", 2)
$d.Content.Find.Execute("64×42=", $true, $false, $false, $false, $false, $true, 1, $false, "79×29=", 2)
$d.Content.Find.Execute("51×12=", $true, $false, $false, $false, $false, $true, 1, $false, "14×58=", 2)
$d.Content.Find.Execute("81×88=", $true, $false, $false, $false, $false, $true, 1, $false, "30×57=", 2)
$d.Content.Find.Execute("91×54=", $true, $false, $false, $false, $false, $true, 1, $false, "62×74=", 2)
$d.Content.Find.Execute("98×90=", $true, $false, $false, $false, $false, $true, 1, $false, "92×32=", 2)
$d.Content.Find.Execute("31×22=", $true, $false, $false, $false, $false, $true, 1, $false, "90×81=", 2)
$d.Content.Find.Execute("70×77=", $true, $false, $false, $false, $false, $true, 1, $false, "68×23=", 2)
$d.Content.Find.Execute("91×74=", $true, $false, $false, $false, $false, $true, 1, $false, "31×20=", 2)
$d.Content.Find.Execute("75×25=", $true, $false, $false, $false, $false, $true, 1, $false, "98×63=", 2)
$d.Content.Find.Execute("69×60=", $true, $false, $false, $false, $false, $true, 1, $false, "25×76=", 2)
$d.Content.Find.Execute("34×58=", $true, $false, $false, $false, $false, $true, 1, $false, "68×83=", 2)
$d.Content.Find.Execute("50×81=", $true, $false, $false, $false, $false, $true, 1, $false, "11×31=", 2)
$d.Content.Find.Execute("20×24=", $true, $false, $false, $false, $false, $true, 1, $false, "56×67=", 2)
$d.Content.Find.Execute("81×54=", $true, $false, $false, $false, $false, $true, 1, $false, "36×97=", 2)
$d.Content.Find.Execute("63×97=", $true, $false, $false, $false, $false, $true, 1, $false, "74×73=", 2)
$d.Content.Find.Execute("35×80=", $true, $false, $false, $false, $false, $true, 1, $false, "82×25=", 2)
$d.Content.Find.Execute("46×39=", $true, $false, $false, $false, $false, $true, 1, $false, "17×90=", 2)
$d.Content.Find.Execute("94×30=", $true, $false, $false, $false, $false, $true, 1, $false, "74×75=", 2)
